# V4 Update: L1 Inductor from 47uF to 100uF (better low current regulation)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the L1 inductor row (row 14) with the new part info:
#   Mfg Part # (E14), Description (F14), LCSC Part # (I14)
$ws.Range("E14").Value = "SWPA5040S101MT"
$ws.Range("F14").Value = "INDUCTOR SMD 100uF 0.75A 5x5mm"
$ws.Range("I14").Value = "C88132"

# Reflect the final selected cell as recorded in the saved workbook
$ws.Range("K14").Select()
